# Regenerate save_data: use K (strikeouts) column values instead of the
# previous "Strike#" derived values. Column G on Sheet1 holds the K values
# for each of the 35 game rows (rows 2-36).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @(1,1,2,3,3,1,0,2,0,2,1,1,1,0,0,1,1,0,1,1,1,3,2,2,2,0,0,4,1,0,2,3,4,2,0)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}

$wb.Save()
